# Updated cryptos list on Wed Jan 17 04:44:28 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.806.12'
$ws.Range('D3').Value = '2.564.90'
$ws.Range('E3').Value = '  +1.39%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.57'
$ws.Range('E5').Value = '  -0.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.19'
$ws.Range('E6').Value = '  +3.37%  '
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.79'
$ws.Range('E10').Value = '  -0.90%  '
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.46'
$ws.Range('E12').Value = '  -1.32%  '
$ws.Range('D13').Value = '2.961.76'
$ws.Range('E13').Value = '  +1.50%  '
$ws.Range('E14').Value = '  -1.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.80'
$ws.Range('E15').Value = '  +3.87%  '
$ws.Range('D16').Value = '2.561.13'
$ws.Range('E16').Value = '  +1.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.846'
$ws.Range('E17').Value = '  -0.78%  '
$ws.Range('D18').Value = '42.842.71'
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.77'
$ws.Range('E19').Value = '  -0.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.54'
$ws.Range('E20').Value = '  -2.65%  '
$ws.Range('E21').Value = '  -0.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.42'
$ws.Range('E22').Value = '  -0.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '249.05'
$ws.Range('E23').Value = '  -1.78%  '
$ws.Range('E24').Value = '  +0.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.06'
$ws.Range('E25').Value = '  -0.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.92'
$ws.Range('E26').Value = '  +0.43%  '
$ws.Range('E27').Value = '  -0.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.39'
$ws.Range('E28').Value = '  -1.86%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '40.15'
$ws.Range('E29').Value = '  -1.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.24'
$ws.Range('E30').Value = '  -0.99%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '157.22'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.79'
$ws.Range('E32').Value = '  -2.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.36'
$ws.Range('E33').Value = '  +0.57%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.11'
$ws.Range('E34').Value = '  -2.83%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0798'
$ws.Range('E35').Value = '  +2.19%  '
$ws.Range('E36').Value = '  +0.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.74'
$ws.Range('E37').Value = '  -1.83%  '
$ws.Range('E38').Value = '  +10.48%  '
$ws.Range('E39').Value = '  -0.13%  '
$ws.Range('E40').Value = '  -0.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '23.34'
$ws.Range('E41').Value = '  -0.35%  '
$ws.Range('E42').Value = '  +7.11%  '
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('E44').Value = '  -0.51%  '
$ws.Range('E45').Value = '  -2.39%  '
$ws.Range('D46').Value = '2.003.53'
$ws.Range('E46').Value = '  -1.60%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.94'
$ws.Range('E47').Value = '  -1.20%  '
$ws.Range('D48').Value = '2.812.62'
$ws.Range('E48').Value = '  +1.35%  '
$ws.Range('E49').Value = '  +2.20%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '74.54'
$ws.Range('E50').Value = '  -0.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '81.61'
$ws.Range('E51').Value = '  -3.99%  '
